$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3028457.8
$ws.Range("I132").Value = 679607.7
$ws.Range("K132").Value = 2038823.1
$ws.Range("M132").Value = -2036293.1
$ws.Range("H141").Value = 5630.1333
$ws.Range("I141").Value = 992.5
$ws.Range("K141").Value = 2977.5
$ws.Range("M141").Value = 2202.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4417.57
$ws.Range("I32").Value = 3776.4421
$ws.Range("J32").Value = 16599
$ws.Range("K32").Value = 3776.4421
$ws.Range("L32").Value = 16599
$ws.Range("M32").Value = -3489.4421
$ws.Range("N32").Value = -17173
$ws.Range("H122").Value = 1670.9584
$ws.Range("I122").Value = 1635.8857
$ws.Range("J122").Value = 1765.3846
$ws.Range("K122").Value = 4907.6571
$ws.Range("L122").Value = 5296.1538
$ws.Range("M122").Value = -2457.6571
$ws.Range("N122").Value = -10196.1538
$ws.Range("H132").Value = 25554294
$ws.Range("I132").Value = 30055368
$ws.Range("K132").Value = 90166104
$ws.Range("M132").Value = -90163574
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 9619439
$ws.Range("I134").Value = 10590447
$ws.Range("J134").Value = 103559.8
$ws.Range("K134").Value = 31771341
$ws.Range("L134").Value = 310679.4
$ws.Range("M134").Value = -31768806
$ws.Range("N134").Value = -315749.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4391677
$ws.Range("I31").Value = 10418243
$ws.Range("J31").Value = 8719.454
$ws.Range("K31").Value = 10418243
$ws.Range("L31").Value = 8719.454
$ws.Range("M31").Value = -10417948
$ws.Range("N31").Value = -9309.454
$ws.Range("H34").Value = 4391677
$ws.Range("I34").Value = 10418243
$ws.Range("J34").Value = 8719.454
$ws.Range("K34").Value = 10418243
$ws.Range("L34").Value = 8719.454
$ws.Range("M34").Value = -10418041
$ws.Range("N34").Value = -9123.454
$ws.Range("H47").Value = 20035.5
$ws.Range("J47").Value = 20035.5
$ws.Range("L47").Value = 20035.5
$ws.Range("N47").Value = -21167.5
$ws.Range("H74").Value = 17132.25
$ws.Range("J74").Value = 18572.908
$ws.Range("L74").Value = 18572.908
$ws.Range("N74").Value = -20320.908
$ws.Range("H77").Value = 17132.25
$ws.Range("J77").Value = 18572.908
$ws.Range("L77").Value = 55718.724
$ws.Range("N77").Value = -64454.724
$ws.Range("H132").Value = 1760.9048
$ws.Range("I132").Value = 1379.3928
$ws.Range("K132").Value = 4138.178400000001
$ws.Range("M132").Value = -1608.178400000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2904.524
$ws.Range("J39").Value = 3136.5264
$ws.Range("L39").Value = 9409.5792
$ws.Range("N39").Value = -9997.5792
$ws.Range("H57").Value = 1005
$ws.Range("I57").Value = 1005
$ws.Range("K57").Value = 3015
$ws.Range("M57").Value = -2456
$ws.Range("H131").Value = 20135.29
$ws.Range("I131").Value = 100282
$ws.Range("J131").Value = 1052.738
$ws.Range("K131").Value = 300846
$ws.Range("L131").Value = 3158.214
$ws.Range("M131").Value = -295806
$ws.Range("N131").Value = -13238.214
$ws.Range("H132").Value = 2027.6086
$ws.Range("I132").Value = 1397.5
$ws.Range("J132").Value = 2160.2632
$ws.Range("K132").Value = 12577.5
$ws.Range("L132").Value = 19442.3688
$ws.Range("M132").Value = -10047.5
$ws.Range("N132").Value = -24502.3688
$ws.Range("H137").Value = 5257.237
$ws.Range("I137").Value = 2261.4285
$ws.Range("J137").Value = 7004.7915
$ws.Range("K137").Value = 6784.2855
$ws.Range("L137").Value = 21014.3745
$ws.Range("M137").Value = -1684.2855
$ws.Range("N137").Value = -31214.3745
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 19817.467
$ws.Range("I113").Value = 2192.7144
$ws.Range("J113").Value = 35239.125
$ws.Range("K113").Value = 2192.7144
$ws.Range("L113").Value = 35239.125
$ws.Range("M113").Value = -22.71439999999984
$ws.Range("N113").Value = -39579.125
$ws.Range("H122").Value = 3243.2974
$ws.Range("I122").Value = 2809.9473
$ws.Range("J122").Value = 3700.7222
$ws.Range("K122").Value = 8429.841899999999
$ws.Range("L122").Value = 11102.1666
$ws.Range("M122").Value = -5979.841899999999
$ws.Range("N122").Value = -16002.1666
$ws.Range("H132").Value = 9958791
$ws.Range("I132").Value = 13034237
$ws.Range("J132").Value = 6063226.5
$ws.Range("K132").Value = 39102711
$ws.Range("L132").Value = 18189679.5
$ws.Range("M132").Value = -39100181
$ws.Range("N132").Value = -18194739.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 16257
$ws.Range("J42").Value = 16257
$ws.Range("L42").Value = 16257
$ws.Range("N42").Value = -17383
$ws.Range("H49").Value = 16257
$ws.Range("J49").Value = 16257
$ws.Range("L49").Value = 16257
$ws.Range("N49").Value = -16551
$ws.Range("H61").Value = 4899.6
$ws.Range("I61").Value = 3249
$ws.Range("J61").Value = 6000
$ws.Range("K61").Value = 3249
$ws.Range("L61").Value = 6000
$ws.Range("M61").Value = -3047
$ws.Range("N61").Value = -6404
$ws.Range("H113").Value = 4899.6
$ws.Range("I113").Value = 3249
$ws.Range("J113").Value = 6000
$ws.Range("K113").Value = 3249
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = -1079
$ws.Range("N113").Value = -10340
$ws.Range("H122").Value = 14753924
$ws.Range("I122").Value = 1522061.1
$ws.Range("J122").Value = 200000000
$ws.Range("K122").Value = 4566183.300000001
$ws.Range("L122").Value = 600000000
$ws.Range("M122").Value = -4563733.300000001
$ws.Range("N122").Value = -600004900
$ws.Range("H132").Value = 6502559.5
$ws.Range("I132").Value = 10998455
$ws.Range("J132").Value = 8488.223
$ws.Range("K132").Value = 32995365
$ws.Range("L132").Value = 25464.669
$ws.Range("M132").Value = -32992835
$ws.Range("N132").Value = -30524.669
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 15625944
$ws.Range("I126").Value = 25000720
$ws.Range("J126").Value = 1316.6666
$ws.Range("K126").Value = 75002160
$ws.Range("L126").Value = 3949.9998
$ws.Range("M126").Value = -74999690
$ws.Range("N126").Value = -8889.9998
$ws.Range("H132").Value = 1326730.6
$ws.Range("I132").Value = 4662.696
$ws.Range("J132").Value = 3665774
$ws.Range("K132").Value = 13988.088
$ws.Range("L132").Value = 10997322
$ws.Range("M132").Value = -11458.088
$ws.Range("N132").Value = -11002382
$ws.Range("H136").Value = 1068.0454
$ws.Range("I136").Value = 1149.8286
$ws.Range("J136").Value = 750
$ws.Range("K136").Value = 3449.4858
$ws.Range("L136").Value = 2250
$ws.Range("M136").Value = -899.4858000000004
$ws.Range("N136").Value = -7350

Write-Host "Applied all changes"